# luban meta-attribute rename:
#   B1 "row:false"    -> "orientation=portrait"
#   C1 "title_rows:4" -> "title_rows=4"
# (everything else on the sheet is unchanged content-wise)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: row:false -> orientation=portrait -------------------------------
$ws.Range("B1").Value = "orientation=portrait"

# Match the font/alignment formatting already used by the other meta cells
# in row 3 (vertically centered, same "等线" font) instead of the plain
# default formatting "row:false" used to have.
$ws.Range("B3").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# The original cell kept the leading "o" in the default run and only
# re-applied the explicit font to the remainder ("rientation=portrait"),
# same pattern Excel produced for the other meta strings on this sheet
# (e.g. "int", "default_avatar", ...).
$ws.Range("B1").Characters(2, 19).Font.Name = "等线"

# --- C1: title_rows:4 -> title_rows=4 -------------------------------------
$ws.Range("C1").Value = "title_rows=4"

# --- Misc view state tweaks recorded in the saved file ---------------------
$ws.Columns(2).ColumnWidth = 19.5
$ws.Range("B13").Select()
